$d = $word.ActiveDocument

# Title cell: "Large Scale " -> "L"  and  "Hate Speech Detection with Cross-Domain Transfer" -> "earning From the Worst (Dynamically generated hate speech dataset)"
$d.Content.Find.Execute("Large Scale ", $true, $false, $false, $false, $false, $true, 1, $false, "L", 2)
$d.Content.Find.Execute("Hate Speech Detection with Cross-Domain Transfer", $true, $false, $false, $false, $false, $true, 1, $false, "earning From the Worst (Dynamically generated hate speech dataset)", 2)

# Size cell
$d.Content.Find.Execute("100k English (27593 hate, 30747 offensive, 41660 none)", $true, $false, $false, $false, $false, $true, 1, $false, "41’225", 2)

# Task Description cell
$d.Content.Find.Execute("Three-class (Hate speech, Offensive language, None)", $true, $false, $false, $false, $false, $true, 1, $false, "Multi-category hate speech detection", 2)
